$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13
$ws.Range("D13").Value = 44383
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 25
$ws.Range("K13").Value = 13000
$ws.Range("L13").Value = 14000
$ws.Range("M13").Value = 13480
$ws.Range("P13").Value = 899

# Row 14
$ws.Range("D14").Value = 44336
$ws.Range("J14").Value = 34
$ws.Range("K14").Value = 24000
$ws.Range("L14").Value = 25000
$ws.Range("M14").Value = 24500
$ws.Range("P14").Value = 1633

# Row 15
$ws.Range("D15").Value = 44341
$ws.Range("J15").Value = 36

# Row 16
$ws.Range("D16").Value = 44400
$ws.Range("J16").Value = 16

# Row 17
$ws.Range("D17").Value = 44442
$ws.Range("J17").Value = 28

# Row 18
$ws.Range("D18").Value = 44453
$ws.Range("J18").Value = 25
$ws.Range("K18").Value = 25000
$ws.Range("L18").Value = 26000
$ws.Range("M18").Value = 25520
$ws.Range("P18").Value = 1701

# Row 19
$ws.Range("D19").Value = 44455
$ws.Range("J19").Value = 18
$ws.Range("K19").Value = 24000
$ws.Range("L19").Value = 25000
$ws.Range("M19").Value = 24500
$ws.Range("P19").Value = 1633
